$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B8").Value = 74
$ws.Range("C8").Value = 7

$ws.Range("C15").Value = 7
$ws.Range("C16").Value = 3

$ws.Range("C20").Value = 5

$ws.Range("B23").Value = 150
$ws.Range("B24").Value = 150
$ws.Range("B25").Value = 150
$ws.Range("B26").Value = 150
$ws.Range("B27").Value = 150
$ws.Range("B28").Value = 150
$ws.Range("B29").Value = 150
$ws.Range("B30").Value = 150
$ws.Range("B31").Value = 150
$ws.Range("B32").Value = 150
$ws.Range("B33").Value = 150
$ws.Range("B34").Value = 150
$ws.Range("B35").Value = 150
$ws.Range("B36").Value = 150

$ws.Range("B37").Value = 100
$ws.Range("B38").Value = 100
$ws.Range("B39").Value = 100
$ws.Range("B40").Value = 100
$ws.Range("B41").Value = 100

$ws.Range("B42").Value = 70

$ws.Range("C22").Select()
